# edit.ps1 - applies the "Vorschlag.docx" revision described by the commit
# "Add files via upload": tweaks to the feature list, a reshuffled list item,
# a renamed section, a reordered "Detailanzeige" field list, a new
# "Statistiken" section, and moving the _GoBack bookmark.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the blank paragraph that sits between
#    "Spiele koennen durch eine Verlinkung per Pfad hinzugefuegt werden."
#    and "Hauptanzeige:" into the following paragraph (i.e. delete the
#    stray empty paragraph mark).
# ---------------------------------------------------------------------
$pHauptanzeige = $d.Content.Find.Execute("Hauptanzeige:")
$blank = $d.Paragraphs.Item(4)
$blank.Range.Delete()

# ---------------------------------------------------------------------
# 2) "Spiel auswaehlen -> evtl. oeffnet Doppelklick die Details Anzeige"
#    becomes "... Doppelklick das Spiel", and a new bullet is inserted
#    right after it: "Spiel auswaehlen -> evtl. rechtsklick
#    detailanzeige/Eigenschaften".
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Spiel auswählen -> evtl. öffnet Doppelklick die Details Anzeige",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Spiel auswählen -> evtl. öffnet Doppelklick das Spiel", 2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Spiel auswählen -> evtl. öffnet Doppelklick das Spiel") {
        $hostIdx = $i
        break
    }
}
$hostPara = $d.Paragraphs.Item($hostIdx)
$hostPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($hostIdx + 1)
$newPara.Range.Text = "Spiel auswählen -> evtl. rechtsklick detailanzeige/Eigenschaften"

# ---------------------------------------------------------------------
# 3) Rename "Detailanzeige" to "Detailanzeige:".
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Detailanzeige", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Detailanzeige:", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Under "Detailanzeige:" the field list is reshuffled:
#      Titel          -> "Titel " (trailing space added)
#      Bild           -> Pfad zur exe
#      Tags           -> Bild
#      Pfad zur exe   -> Tags
#    Editing Range.Text in place keeps each paragraph's own run
#    formatting (e.g. the yellow highlight) exactly where it was.
# ---------------------------------------------------------------------
$titelFound = $d.Content.Find.Execute("Detailanzeige:")
$searchRange = $d.Content
$searchRange.Start = $d.Content.Find.Execute("Detailanzeige:") 

# locate the four paragraphs following "Detailanzeige:" directly by index
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Detailanzeige:") {
        $detailIdx = $i
        break
    }
}

$pTitel = $d.Paragraphs.Item($detailIdx + 1)
$pBild  = $d.Paragraphs.Item($detailIdx + 2)
$pTags  = $d.Paragraphs.Item($detailIdx + 3)
$pPfad  = $d.Paragraphs.Item($detailIdx + 4)

$pTitel.Range.Text = "Titel "
$pBild.Range.Text  = "Pfad zur exe"
$pTags.Range.Text  = "Bild"
$pPfad.Range.Text  = "Tags"

# ---------------------------------------------------------------------
# 5) Move the _GoBack bookmark from the trailing empty paragraph to the
#    end of "Features:" (it now sits at the top of the document).
# ---------------------------------------------------------------------
$featuresRange = $d.Content
$featuresRange.Find.Execute("Features:") | Out-Null
$featuresRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $featuresRange) | Out-Null

# ---------------------------------------------------------------------
# 6) Append the new "Statistiken:" section with its own bullet list at
#    the end of the document.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$statistikenPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$statistikenPara.Range.Text = "Statistiken:"

$statistikenPara.Range.InsertParagraphAfter()
$bulletPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bulletPara.Range.Text = "Spielstunden (Letze zwei Wochen?)"
$bulletPara.Style = "Listenabsatz"
$sourceList = $d.Paragraphs.Item($detailIdx - 1).Range.ListFormat.ListTemplate
$bulletPara.Range.ListFormat.ApplyListTemplateWithLevel($sourceList)
